$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'2026-02-28 06:48:34"
$ws.Range("E3").Value = "'2026-02-28 06:48:37"
$ws.Range("H3").Value = "'86%"
$ws.Range("O3").Value = "'-0.7 °C"
$ws.Range("E4").Value = "'2026-02-28 06:48:40"
$ws.Range("H4").Value = "'94%"
$ws.Range("K4").Value = "'-0.1 MJ/m2"
$ws.Range("O4").Value = "'8.0 °C"
$ws.Range("E5").Value = "'2026-02-28 06:48:42"
$ws.Range("H5").Value = "'97%"
$ws.Range("N5").Value = "'-1.7 °C 6:27 TU"
$ws.Range("O5").Value = "'-0.3 °C"
$ws.Range("E6").Value = "'2026-02-28 06:48:45"
$ws.Range("N6").Value = "'9.9 °C 6:23 TU"
$ws.Range("E7").Value = "'2026-02-28 06:48:48"
$ws.Range("J7").Value = "'1023.8 hPa"
$ws.Range("E8").Value = "'2026-02-28 06:48:50"
$ws.Range("M8").Value = "'9.0 °C 6:20 TU"
$ws.Range("E9").Value = "'2026-02-28 06:48:53"
$ws.Range("O9").Value = "'7.2 °C"
$ws.Range("E10").Value = "'2026-02-28 06:48:55"
$ws.Range("L10").Value = "'6.5 km/h - 46º 6:12 TU"
$ws.Range("M10").Value = "'9.9 °C 6:22 TU"
$ws.Range("O10").Value = "'8.0 °C"
$ws.Range("E11").Value = "'2026-02-28 06:48:58"
$ws.Range("N11").Value = "'0.6 °C 6:17 TU"
$ws.Range("O11").Value = "'3.0 °C"
$ws.Range("E12").Value = "'2026-02-28 06:49:00"
$ws.Range("E13").Value = "'2026-02-28 06:49:03"
$ws.Range("J13").Value = "'1026.5 hPa"
$ws.Range("N13").Value = "'-1.0 °C 6:03 TU"
$ws.Range("O13").Value = "'0.9 °C"
$ws.Range("E14").Value = "'2026-02-28 06:49:05"
$ws.Range("H14").Value = "'96%"
$ws.Range("L14").Value = "'22.0 km/h - 95º 6:18 TU"
$ws.Range("M14").Value = "'12.7 °C 6:29 TU"
$ws.Range("O14").Value = "'10.7 °C"
$ws.Range("E15").Value = "'2026-02-28 06:49:08"
$ws.Range("E16").Value = "'2026-02-28 06:49:10"
$ws.Range("H16").Value = "'58%"
$ws.Range("E17").Value = "'2026-02-28 06:49:13"
$ws.Range("H17").Value = "'47%"
$ws.Range("N17").Value = "'2.0 °C 6:28 TU"
$ws.Range("O17").Value = "'4.3 °C"
$ws.Range("E18").Value = "'2026-02-28 06:49:16"
$ws.Range("L18").Value = "'4.3 km/h - 263º 6:10 TU"
$ws.Range("O18").Value = "'8.5 °C"
$ws.Range("E19").Value = "'2026-02-28 06:49:18"
$ws.Range("H19").Value = "'68%"
$ws.Range("N19").Value = "'4.3 °C 6:18 TU"
$ws.Range("O19").Value = "'7.7 °C"
$ws.Range("E20").Value = "'2026-02-28 06:49:21"
$ws.Range("H20").Value = "'37%"
$ws.Range("L20").Value = "'23.0 km/h - 159º 6:28 TU"
$ws.Range("N20").Value = "'-1.5 °C 6:07 TU"
$ws.Range("O20").Value = "'-0.1 °C"
$ws.Range("E21").Value = "'2026-02-28 06:49:24"
$ws.Range("J21").Value = "'1024.3 hPa"
$ws.Range("O21").Value = "'4.8 °C"
$ws.Range("E22").Value = "'2026-02-28 06:49:26"
$ws.Range("H22").Value = "'57%"
$ws.Range("N22").Value = "'-1.9 °C 6:29 TU"
$ws.Range("O22").Value = "'-0.9 °C"
$ws.Range("E23").Value = "'2026-02-28 06:49:29"
$ws.Range("H23").Value = "'67%"
$ws.Range("N23").Value = "'-0.9 °C 6:29 TU"
$ws.Range("O23").Value = "'-0.1 °C"
$ws.Range("E24").Value = "'2026-02-28 06:49:31"
$ws.Range("J24").Value = "'1023.7 hPa"
$ws.Range("O24").Value = "'6.3 °C"
$ws.Range("E25").Value = "'2026-02-28 06:49:34"
$ws.Range("H25").Value = "'54%"
$ws.Range("O25").Value = "'0.6 °C"
$ws.Range("E26").Value = "'2026-02-28 06:49:37"
$ws.Range("O26").Value = "'4.3 °C"
$ws.Range("E27").Value = "'2026-02-28 06:49:39"
$ws.Range("O27").Value = "'2.0 °C"
$ws.Range("E28").Value = "'2026-02-28 06:49:42"
$ws.Range("J28").Value = "'1024.7 hPa"
$ws.Range("O28").Value = "'6.4 °C"
$ws.Range("E29").Value = "'2026-02-28 06:49:44"
$ws.Range("E30").Value = "'2026-02-28 06:49:47"
$ws.Range("J30").Value = "'1024.4 hPa"
$ws.Range("E31").Value = "'2026-02-28 06:49:50"
$ws.Range("N31").Value = "'9.5 °C 6:29 TU"
$ws.Range("E32").Value = "'2026-02-28 06:49:52"
$ws.Range("H32").Value = "'92%"
$ws.Range("E33").Value = "'2026-02-28 06:49:55"
$ws.Range("H33").Value = "'72%"
$ws.Range("J33").Value = "'1023.8 hPa"
$ws.Range("N33").Value = "'3.2 °C 6:11 TU"
$ws.Range("O33").Value = "'4.8 °C"
$ws.Range("E34").Value = "'2026-02-28 06:49:57"
$ws.Range("H34").Value = "'64%"
$ws.Range("N34").Value = "'-2.0 °C 6:15 TU"
$ws.Range("O34").Value = "'-0.3 °C"
$ws.Range("E35").Value = "'2026-02-28 06:50:00"
$ws.Range("H35").Value = "'84%"
$ws.Range("J35").Value = "'1023.1 hPa"
$ws.Range("N35").Value = "'4.7 °C 6:29 TU"
$ws.Range("O35").Value = "'6.4 °C"
$ws.Range("E36").Value = "'2026-02-28 06:50:03"
$ws.Range("J36").Value = "'1024.4 hPa"
$ws.Range("O36").Value = "'9.8 °C"
$ws.Range("E37").Value = "'2026-02-28 06:50:06"
$ws.Range("J37").Value = "'1026.1 hPa"
$ws.Range("N37").Value = "'3.0 °C 6:29 TU"
$ws.Range("O37").Value = "'4.4 °C"
$ws.Range("E38").Value = "'2026-02-28 06:50:08"
$ws.Range("O38").Value = "'9.2 °C"
$ws.Range("E39").Value = "'2026-02-28 06:50:10"
$ws.Range("H39").Value = "'46%"
$ws.Range("E40").Value = "'2026-02-28 06:50:13"
$ws.Range("J40").Value = "'1024.8 hPa"
$ws.Range("O40").Value = "'3.2 °C"
$ws.Range("E41").Value = "'2026-02-28 06:50:16"
$ws.Range("M41").Value = "'12.1 °C 6:29 TU"
$ws.Range("E42").Value = "'2026-02-28 06:50:18"
$ws.Range("E43").Value = "'2026-02-28 06:50:21"
$ws.Range("H43").Value = "'87%"
$ws.Range("N43").Value = "'1.8 °C 6:20 TU"
$ws.Range("O43").Value = "'3.6 °C"
$ws.Range("E44").Value = "'2026-02-28 06:50:23"
$ws.Range("E45").Value = "'2026-02-28 06:50:25"
$ws.Range("G45").Value = "'2 cm"
$ws.Range("J45").Value = "'1024.2 hPa"
$ws.Range("N45").Value = "'4.8 °C 6:12 TU"
$ws.Range("O45").Value = "'6.5 °C"
$ws.Range("E46").Value = "'2026-02-28 06:50:28"
$ws.Range("J46").Value = "'1023.3 hPa"
$ws.Range("O46").Value = "'11.0 °C"
